$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks like a plain number but must remain text
# (matching the source data which stores these as inline strings, not numbers).
# We briefly mark the cell as Text-formatted, write the string, then restore the
# default (unstyled) cell style so the saved style table stays clean.

$ws.Range("D2").Value = "26.976.36"
$ws.Range("E2").Value = "  -0.20%  "

$ws.Range("D3").Value = "1.676.90"
$ws.Range("E3").Value = "  +0.17%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.21"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.41%  "

$ws.Range("E6").Value = "  +1.47%  "

$ws.Range("E7").Value = "  +0.00%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.252"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.22%  "

$ws.Range("E9").Value = "  +0.38%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.32"
$ws.Range("D10").Style = "Normal"

$ws.Range("E11").Value = "  +0.07%  "

$ws.Range("D12").Value = "1.914.47"
$ws.Range("E12").Value = "  +0.26%  "

$ws.Range("D13").Value = "1.685.39"
$ws.Range("E13").Value = "  +0.66%  "

$ws.Range("E14").Value = "  -0.06%  "

$ws.Range("E15").Value = "  +1.31%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.87"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.02%  "

$ws.Range("D17").Value = "26.992.49"
$ws.Range("E17").Value = "  -0.14%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "237.11"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.44%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "8.09"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.55%  "

$ws.Range("D20").Value = "0.0₃0733"
$ws.Range("E20").Value = "  -0.63%  "

$ws.Range("E21").Value = "  -0.03%  "

$ws.Range("E22").Value = "  -0.73%  "

$ws.Range("E23").Value = "  -0.94%  "

$ws.Range("E24").Value = "  -1.54%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "145.63"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.03%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.24"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.16%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.10"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.98%  "

$ws.Range("E28").Value = "  -1.26%  "

$ws.Range("E29").Value = "  -0.02%  "

$ws.Range("E30").Value = "  -0.15%  "

$ws.Range("E31").Value = "  -0.47%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.32"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.21%  "

$ws.Range("D33").Value = "1.486.07"
$ws.Range("E33").Value = "  +0.81%  "

$ws.Range("E34").Value = "  +0.72%  "

$ws.Range("E35").Value = "  +4.74%  "

$ws.Range("E36").Value = "  +0.07%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.586"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.03%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0175"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.02%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.903"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.60%  "

$ws.Range("E40").Value = "  -3.99%  "

$ws.Range("E41").Value = "  +1.09%  "

$ws.Range("E42").Value = "  -0.02%  "

$ws.Range("E43").Value = "  +1.83%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "67.51"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.45%  "

$ws.Range("D45").Value = "1.820.40"
$ws.Range("E45").Value = "  -0.01%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.780"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.10%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "90.63"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.51%  "

$ws.Range("E48").Value = "  +2.25%  "

$ws.Range("E49").Value = "  -0.39%  "

$ws.Range("E50").Value = "  +1.91%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0509"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.38%  "
